# "added functionality to pick best lot"
# The day/time column headers in row 1 (D1:W1) are renamed from the
# "DayHH0" form (e.g. "Mon08", "Mon010") to a clearer "Day-H" form
# (e.g. "Mon-8", "Mon-10") so the lot-picking logic downstream can key off
# a consistent "<Day>-<Hour>" label.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D1").Value = "Mon-8"
$ws.Range("E1").Value = "Mon-10"
$ws.Range("F1").Value = "Mon-12"
$ws.Range("G1").Value = "Mon-2"

$ws.Range("H1").Value = "Tue-8"
$ws.Range("I1").Value = "Tue-10"
$ws.Range("J1").Value = "Tue-12"
$ws.Range("K1").Value = "Tue-2"

$ws.Range("L1").Value = "Wed-8"
$ws.Range("M1").Value = "Wed-10"
$ws.Range("N1").Value = "Wed-12"
$ws.Range("O1").Value = "Wed-2"

$ws.Range("P1").Value = "Thu-8"
$ws.Range("Q1").Value = "Thu-10"
$ws.Range("R1").Value = "Thu-12"
$ws.Range("S1").Value = "Thu-2"

$ws.Range("T1").Value = "Fri-8"
$ws.Range("U1").Value = "Fri-10"
$ws.Range("V1").Value = "Fri-12"
$ws.Range("W1").Value = "Fri-2"

# Move the sheet's selection to where the author's cursor ended up
# (U19) after reviewing/editing the new headers.
$ws.Range("U19").Select()
